# Update trading results log - append two new rows (100, 101) describing a
# new TRADING_ATTEMPT / POSITION_FAILED pair for ETH, matching the existing
# log row layout (timestamp, action, token, signal_type, price,
# position_size_usd, leverage, stiffness, pnl_percent, exit_reason, status,
# error_message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A lone apostrophe is Excel's "treat as text" quote-prefix; assigning it as
# a cell's Value yields a genuine empty-string text cell (matching the blank
# columns already present on every other log row) instead of clearing the
# cell outright the way Value = "" would.
$apos = [char]39

# Row 100: TRADING_ATTEMPT
$ws.Range("A100").Value = "2025-10-28T12:44:30.584659"
$ws.Range("B100").Value = "TRADING_ATTEMPT"
$ws.Range("C100").Value = "ETH"
$ws.Range("D100").Value = "UNKNOWN"
$ws.Range("E100").Value = 4123.206220671679
$ws.Range("F100").Value = $apos
$ws.Range("G100").Value = $apos
$ws.Range("H100").Value = $apos
$ws.Range("I100").Value = $apos
$ws.Range("J100").Value = $apos
$ws.Range("K100").Value = "ATTEMPT"
$ws.Range("L100").Value = "Attempting trade 1/1"

# Row 101: POSITION_FAILED
$ws.Range("A101").Value = "2025-10-28T12:44:32.730295"
$ws.Range("B101").Value = "POSITION_FAILED"
$ws.Range("C101").Value = "ETH"
$ws.Range("D101").Value = "UNKNOWN"
$ws.Range("E101").Value = $apos
$ws.Range("F101").Value = $apos
$ws.Range("G101").Value = $apos
$ws.Range("H101").Value = $apos
$ws.Range("I101").Value = $apos
$ws.Range("J101").Value = $apos
$ws.Range("K101").Value = "FAILED"
$ws.Range("L101").Value = "Trade execution failed for trade 1"

# Drop the quote-prefix formatting the apostrophe trick applied so the new
# blank cells carry the default (unstyled) format, same as the rest of the
# sheet's data rows.
$ws.Range("F100:J100").ClearFormats()
$ws.Range("E101:J101").ClearFormats()
